$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 175
$ws1.Range("F3").Value = 1335
$ws1.Range("F5").Value = 934
$ws1.Range("F7").Value = 216
$ws1.Range("F8").Value = 551
$ws1.Range("F12").Value = 3093
$ws1.Range("F13").Value = 2703
$ws1.Range("F20").Value = 5532
$ws1.Range("F21").Value = 602
$ws1.Range("F25").Value = 417
$ws1.Range("F26").Value = 1185
$ws1.Range("F28").Value = 98

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 10
$ws2.Range("F9").Value = 46
$ws2.Range("F25").Value = 4021

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value = 1405
$ws3.Range("F10").Value = 394

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1405
$ws4.Range("F8").Value = 394
$ws4.Range("F10").Value = 175
$ws4.Range("F11").Value = 1335
$ws4.Range("F12").Value = 934
$ws4.Range("F15").Value = 216
$ws4.Range("F16").Value = 551
$ws4.Range("F18").Value = 3093
$ws4.Range("F19").Value = 2703
$ws4.Range("F23").Value = 46
$ws4.Range("F26").Value = 5532
$ws4.Range("F28").Value = 602
$ws4.Range("F33").Value = 417
$ws4.Range("F40").Value = 1185
$ws4.Range("F48").Value = 98
